# Weekly driver report update for 2025-04-28
# Updates the "Bad Drivers" and "Good Drivers" tables on the single
# "Driver Summary" worksheet: new adapter rows are inserted into both
# tables, and sample counts / percentages / vintages are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room for the new rows.
#    "Bad Drivers" table gains 2 new adapter rows (rows 8 & 9) before the
#    existing "Totals:" row, which is currently row 8.
# ---------------------------------------------------------------------
$ws.Range("A8:A9").EntireRow.Insert()

#    "Good Drivers" table gains 1 new adapter row. After the insert above,
#    its first data row ("...23.100.0.4", formerly row 16) now sits at
#    row 18 - insert a fresh row there, ahead of it.
$ws.Range("A18").EntireRow.Insert()

# ---------------------------------------------------------------------
# 2) "Bad Drivers" table (header at row 2, data rows 3-9, totals row 10)
# ---------------------------------------------------------------------
$ws.Range("C3").Value = 861
$ws.Range("D3").Value = 78.09999999999999

$ws.Range("A4").Value = "Killer Wi-Fi 6 AX500-DBS Wireless Network Adapter - 1.0.0.1769"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 1084
$ws.Range("D4").Value = 91.8

$ws.Range("A5").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.80.0.7"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 772
$ws.Range("D5").Value = 96

$ws.Range("A6").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.0.10"
$ws.Range("B6").Value = 11
$ws.Range("C6").Value = 3382
$ws.Range("D6").Value = 96.40000000000001

$ws.Range("A7").Value = "NETGEAR A8000 WiFi 6 & 6E Adapter - 1.0.0.108"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 86
$ws.Range("D7").Value = 97.59999999999999

$ws.Range("A8").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.250.10.1"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 36
$ws.Range("D8").Value = 97.90000000000001

$ws.Range("A9").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.40.0.4"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 37
$ws.Range("D9").Value = 98.7

$ws.Range("B10").Value = 24
$ws.Range("C10").Value = 6258

# ---------------------------------------------------------------------
# 3) "Good Drivers" table (header at row 17, data rows 18-27)
#
#    The "Driver Vintage" column (E) stores dates as plain text
#    (e.g. "2024-11-10"), not real Excel date serials. Writing a
#    date-shaped string straight into `.Value` gets auto-converted to a
#    date by Excel, so each vintage cell is briefly marked as Text
#    (NumberFormat "@") while the literal string is assigned, then
#    ClearFormats()+re-align restores the plain right-aligned look used
#    by the rest of the column without leaving the cell tagged as Text.
# ---------------------------------------------------------------------
function Set-VintageText($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
    $cell.HorizontalAlignment = -4152
}

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B18").Value = 11128
$ws.Range("B18").NumberFormat = "#,##0"
$ws.Range("B18").HorizontalAlignment = -4152
$ws.Range("D18").Value = 100

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B19").Value = 486214
$ws.Range("D19").Value = 99.90000000000001
Set-VintageText $ws.Range("E19") "2024-11-10"

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B20").Value = 11140
$ws.Range("D20").Value = 100
Set-VintageText $ws.Range("E20") "2022-08-29"

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B21").Value = 14487
$ws.Range("D21").Value = 100
Set-VintageText $ws.Range("E21") "2022-05-23"

$ws.Range("A22").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B22").Value = 265400
$ws.Range("D22").Value = 99.90000000000001
Set-VintageText $ws.Range("E22") "2022-05-01"

$ws.Range("A23").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B23").Value = 79953
$ws.Range("D23").Value = 99.90000000000001
Set-VintageText $ws.Range("E23") "2021-08-18"

$ws.Range("A24").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B24").Value = 35355
$ws.Range("D24").Value = 100
Set-VintageText $ws.Range("E24") "2021-04-27"

$ws.Range("A25").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B25").Value = 65425
$ws.Range("D25").Value = 100
Set-VintageText $ws.Range("E25") "2020-08-05"

$ws.Range("A26").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B26").Value = 117653
$ws.Range("D26").Value = 100
Set-VintageText $ws.Range("E26") "2020-01-06"

$ws.Range("A27").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B27").Value = 56018
$ws.Range("D27").Value = 100
Set-VintageText $ws.Range("E27") "2019-12-14"

Write-Host "Driver summary refreshed for 2025-04-28"
